$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "와이바이오로직스"
$ws.Range("B2").Value = "2023.11.10~11.16"
$ws.Range("C2").Value = "9,000~11,000"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 13500
$ws.Range("F2").Value = "유안타증권"

$ws.Range("A3").Value = "한선엔지니어링"
$ws.Range("B3").Value = "2023.11.02~11.08"
$ws.Range("C3").Value = "5,200~6,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 22100
$ws.Range("F3").Value = "대신증권"

$ws.Range("A4").Value = "에코아이"
$ws.Range("B4").Value = "2023.11.01~11.07"
$ws.Range("C4").Value = "28,500~34,700"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 59251
$ws.Range("F4").Value = "KB증권"

$ws.Range("A5").Value = "동인기연(유가)"
$ws.Range("B5").Value = "2023.11.01~11.07"
$ws.Range("C5").Value = "33,000~37,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 60654
$ws.Range("F5").Value = "NH투자증권"

$ws.Range("A6").Value = "스톰테크"
$ws.Range("B6").Value = "2023.10.31~11.06"
$ws.Range("C6").Value = "8,000~9,500"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 26800
$ws.Range("F6").Value = "하이투자증권"

$ws.Range("A7").Value = "블루엠텍"
$ws.Range("B7").Value = "2023.10.31~11.06"
$ws.Range("C7").Value = "15,000~19,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 21000
$ws.Range("F7").Value = "하나증권,키움증권"

$ws.Range("A8").Value = "에코프로머티리얼즈"
$ws.Range("B8").Value = "2023.10.30~11.03"
$ws.Range("C8").Value = "36,200~44,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 524031
$ws.Range("F8").Value = "미래에셋증권,NH투자증권,하이투자증권"

$ws.Range("A9").Value = "캡스톤파트너스"
$ws.Range("B9").Value = "2023.10.26~11.01"
$ws.Range("C9").Value = "3,200~3,600"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = 5107
$ws.Range("F9").Value = "NH투자증권"

$ws.Range("A10").Value = "에이텀"
$ws.Range("B10").Value = "2023.10.26~11.01"
$ws.Range("C10").Value = "23,000~30,000"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = 14950
$ws.Range("F10").Value = "하나증권"

$ws.Range("A11").Value = "한국스팩13호"
$ws.Range("B11").Value = "2023.10.25~10.26"
$ws.Range("C11").Value = "2,000~2,000"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = 8000
$ws.Range("F11").Value = "한국투자증권"

$ws.Range("A12").Value = "그린리소스"
$ws.Range("B12").Value = "2023.10.25~10.31"
$ws.Range("C12").Value = "11,000~14,000"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = 18040
$ws.Range("F12").Value = "NH투자증권"

$ws.Range("A13").Value = "에이직랜드"
$ws.Range("B13").Value = "2023.10.23~10.27"
$ws.Range("C13").Value = "19,100~21,400"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 50353
$ws.Range("F13").Value = "삼성증권"

$ws.Range("A14").Value = "에스와이스틸텍"
$ws.Range("B14").Value = "2023.10.23~10.27"
$ws.Range("C14").Value = "1,200~1,500"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = 8400
$ws.Range("F14").Value = "KB증권"

$ws.Range("A15").Value = "컨텍"
$ws.Range("B15").Value = "2023.10.20~10.26"
$ws.Range("C15").Value = "20,300~22,500"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = 41818
$ws.Range("F15").Value = "대신증권"

$ws.Range("A16").Value = "큐로셀"
$ws.Range("B16").Value = "2023.10.20~10.26"
$ws.Range("C16").Value = "29,800~33,500"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = 47680
$ws.Range("F16").Value = "미래에셋증권,삼성증권"

$ws.Range("A17").Value = "메가터치"
$ws.Range("B17").Value = "2023.10.20~10.26"
$ws.Range("C17").Value = "3,500~4,000"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = 18200
$ws.Range("F17").Value = "NH투자증권"

$ws.Range("A18").Value = "비아이매트릭스"
$ws.Range("B18").Value = "2023.10.19~10.25"
$ws.Range("C18").Value = "9,100~11,000"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = 10920
$ws.Range("F18").Value = "IBK투자증권"

$ws.Range("A19").Value = "KB스팩27호"
$ws.Range("B19").Value = "2023.10.19~10.20"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = 25000
$ws.Range("F19").Value = "KB증권"

$ws.Range("A20").Value = "유투바이오"
$ws.Range("B20").Value = "2023.10.18~10.19"
$ws.Range("C20").Value = "3,300~3,900"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = 3724
$ws.Range("F20").Value = "신한투자증권"

$ws.Range("A21").Value = "쏘닉스"
$ws.Range("B21").Value = "2023.10.17~10.23"
$ws.Range("C21").Value = "5,000~7,000"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = 18000
$ws.Range("F21").Value = "KB증권"
